$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet tab: "Apr 14 2022" -> "Apr 16 2022"
$ws.Name = "Apr 16 2022"

# 2. Update the plan date (A2) two days forward (Apr 14 -> Apr 16, 2022)
$ws.Range("A2").Value = 44667

# 3. Row 4: clear the "Notes" cell (E4) entirely
$ws.Range("E4").Clear()

# 4. Row 4: "Tested by" (F4) N/A -> Nick
$ws.Range("F4").Value = "Nick"

# 5. Row 4: results G4:I4 Failed -> Passed, with a new "good" look (green fill)
$results = $ws.Range("G4:I4")
$results.Value = "Passed"
$results.Interior.ThemeColor = 10
$results.HorizontalAlignment = -4108
$results.VerticalAlignment = -4108
$results.WrapText = $true

# 6. Update the active selection to I9
$ws.Range("I9").Select()
